$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-21) down by 3 rows to make room for 3 new rows
# at the top. Work from the bottom up so we never overwrite a row before reading it.
for ($i = 21; $i -ge 2; $i--) {
    $dest = $i + 3
    $a = $ws.Range("A$i").Value()
    $b = $ws.Range("B$i").Value()
    $c = $ws.Range("C$i").Value()
    $ws.Range("A$dest").Value = $a
    $ws.Range("B$dest").Value = $b
    $ws.Range("C$dest").Value = $c
}

# New rows of data inserted at the top (new rows 2-4)
$topData = @(
    @(-0.0074830991216003, 0.0424551330506801, 0.0247400421649217),
    @(0.0018325957935303, 0.0197004042565822, 0.0239764600992202),
    @(0.0056505035609006, -0.0122173046693205, 0.009315694682300001)
)

$r = 2
foreach ($row in $topData) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r = $r + 1
}

# New rows of data appended at the bottom (rows 25-31)
$bottomData = @(
    @(0.0568104684352874, 0.1012509167194366, -0.0526871271431446),
    @(0.0032070425804704, -0.0478002056479454, -0.024892758578062),
    @(-0.0117591563612222, 0.0113010071218013, 0.0294742472469806),
    @(-0.0128281703218817, -0.0500909499824047, -0.0126754539087414),
    @(0.0056505035609006, -0.0193949714303016, -0.0198531206697225),
    @(0.0088575463742017, 0.0649044290184974, 0.0067195175215601),
    @(-0.0021380283869802, 0.01328631862998, 0.0001527163112768)
)

$r = 25
foreach ($row in $bottomData) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r = $r + 1
}
